$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers (e.g. "57.72", "1.00")
# need an explicit Text number format while the value is assigned, otherwise
# Excel auto-converts them to doubles and loses the original text formatting
# (trailing zeros, etc). The format is reset back to Normal afterwards so the
# cell style matches the original (unstyled) cells.
$textCells = @(
    "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D16", "D17", "D20", "D21", "D23", "D25", "D27", "D28", "D31", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D50"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "37.426.32"
$ws.Range("E2").Value = "  +3.70%  "

# Row 3
$ws.Range("D3").Value = "2.065.88"
$ws.Range("E3").Value = "  +6.03%  "

# Row 4
$ws.Range("E4").Value = "  -0.17%  "

# Row 5
$ws.Range("D5").Value = "236.22"
$ws.Range("E5").Value = "  +3.59%  "

# Row 6
$ws.Range("D6").Value = "0.616"
$ws.Range("E6").Value = "  +4.56%  "

# Row 7
$ws.Range("B7").Value = "Solana"
$ws.Range("C7").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D7").Value = "57.72"
$ws.Range("E7").Value = "  +9.53%  "

# Row 8
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.06%  "

# Row 9
$ws.Range("D9").Value = "0.381"
$ws.Range("E9").Value = "  +5.11%  "

# Row 10
$ws.Range("D10").Value = "57.77"
$ws.Range("E10").Value = "  +2.16%  "

# Row 11
$ws.Range("D11").Value = "0.0760"
$ws.Range("E11").Value = "  +4.52%  "

# Row 12
$ws.Range("D12").Value = "0.102"
$ws.Range("E12").Value = "  +4.85%  "

# Row 13
$ws.Range("D13").Value = "2.368.83"
$ws.Range("E13").Value = "  +5.95%  "

# Row 14
$ws.Range("D14").Value = "14.31"
$ws.Range("E14").Value = "  +4.63%  "

# Row 15
$ws.Range("D15").Value = "20.85"
$ws.Range("E15").Value = "  +8.34%  "

# Row 16
$ws.Range("D16").Value = "0.775"
$ws.Range("E16").Value = "  +5.25%  "

# Row 17
$ws.Range("D17").Value = "5.18"
$ws.Range("E17").Value = "  +5.21%  "

# Row 18
$ws.Range("D18").Value = "2.067.87"
$ws.Range("E18").Value = "  +5.58%  "

# Row 19
$ws.Range("D19").Value = "37.554.07"
$ws.Range("E19").Value = "  +4.13%  "

# Row 20
$ws.Range("D20").Value = "6.14"
$ws.Range("E20").Value = "  +25.04%  "

# Row 21
$ws.Range("D21").Value = "68.51"
$ws.Range("E21").Value = "  +2.92%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0810"
$ws.Range("E22").Value = "  +3.63%  "

# Row 23
$ws.Range("D23").Value = "225.06"
$ws.Range("E23").Value = "  +2.91%  "

# Row 24
$ws.Range("E24").Value = "  -0.06%  "

# Row 25
$ws.Range("D25").Value = "2.45"
$ws.Range("E25").Value = "  +7.19%  "

# Row 26
$ws.Range("E26").Value = "  +3.01%  "

# Row 27
$ws.Range("D27").Value = "163.50"
$ws.Range("E27").Value = "  +2.34%  "

# Row 28
$ws.Range("D28").Value = "8.83"
$ws.Range("E28").Value = "  +5.86%  "

# Row 29
$ws.Range("E29").Value = "  +11.41%  "

# Row 30
$ws.Range("E30").Value = "  +10.02%  "

# Row 31
$ws.Range("D31").Value = "19.19"
$ws.Range("E31").Value = "  +3.43%  "

# Row 32
$ws.Range("E32").Value = "  +2.59%  "

# Row 33
$ws.Range("E33").Value = "  +18.26%  "

# Row 34
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "4.46"
$ws.Range("E34").Value = "  +4.91%  "

# Row 35
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "0.0626"
$ws.Range("E35").Value = "  +5.28%  "

# Row 36
$ws.Range("D36").Value = "4.46"
$ws.Range("E36").Value = "  +7.69%  "

# Row 37
$ws.Range("B37").Value = "BinanceUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  -0.04%  "

# Row 38
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").Value = "1.80"
$ws.Range("E38").Value = "  +1.10%  "

# Row 39
$ws.Range("D39").Value = "3.35"
$ws.Range("E39").Value = "  +9.21%  "

# Row 40
$ws.Range("D40").Value = "5.84"
$ws.Range("E40").Value = "  +17.69%  "

# Row 41
$ws.Range("D41").Value = "2.98"
$ws.Range("E41").Value = "  -0.24%  "

# Row 42
$ws.Range("B42").Value = "FTXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D42").Value = "4.46"
$ws.Range("E42").Value = "  +31.76%  "

# Row 43
$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").Value = "0.0958"
$ws.Range("E43").Value = "  +12.45%  "

# Row 44
$ws.Range("D44").Value = "1.466.83"
$ws.Range("E44").Value = "  +5.36%  "

# Row 45
$ws.Range("D45").Value = "95.43"
$ws.Range("E45").Value = "  +11.76%  "

# Row 46
$ws.Range("D46").Value = "0.0211"
$ws.Range("E46").Value = "  +7.55%  "

# Row 47
$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D47").Value = "1.13"
$ws.Range("E47").Value = "  +6.95%  "

# Row 48
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "16.09"
$ws.Range("E48").Value = "  +11.42%  "

# Row 49
$ws.Range("E49").Value = "  +5.93%  "

# Row 50
$ws.Range("D50").Value = "7.27"
$ws.Range("E50").Value = "  +9.11%  "

# Row 51
$ws.Range("E51").Value = "  +2.74%  "

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
